$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new project row (row 4: "FussionFall Re:Spawn") to the Projects sheet.
# Row 4 already exists in the sheet (as a mostly-empty placeholder row with
# only A4:C4 present), so we fill in all 8 columns (A-H) with the new
# project's data, matching the layout of rows 2 and 3.
# ---------------------------------------------------------------------------

function Copy-FormatOnly {
    param($srcAddr, $dstAddr)
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# Pre-seed the formatting for row 4 from row 3 (same visual layout: A/B/C/D/E/F
# centered+wrapped, G centered+wrapped hyperlink style) before writing values,
# so the existing style table entries get reused instead of new ones created.
Copy-FormatOnly "A3" "A4"
Copy-FormatOnly "B3" "B4"
Copy-FormatOnly "C3" "C4"
Copy-FormatOnly "D3" "D4"
Copy-FormatOnly "E3" "E4"
Copy-FormatOnly "F3" "F4"
Copy-FormatOnly "G3" "G4"

$ws.Range("A4").Value = "/img/ffrespawn.png"
$ws.Range("B4").Value = "In game screenshot of FusionFall Re:Spawn"
$ws.Range("C4").Value = "FussionFall Re:Spawn"
$ws.Range("D4").Value = "game dev, c#, unity"
$ws.Range("E4").Value = "FusionFall Re:Spawn was a fan-made revival project of the online MMO FusionFall."
$ws.Range("F4").Value = "<p>The aim of the project has been to recreate the original game with a more recent and stable version of Unity. Rewriting the story while loosely following the original to fill narrative gaps and to leave room for exciting new stories to take place in the universe</p><p>The project was being developed by volunteers whose vision was to bring the game they fell in love to the newer generations. I joined the team for a similar reason: I joined the development team in hope of creating something that people will be able to enjoy as much as I enjoyed playing the original game as a kid back in 2012-2013. Sadly, the project was discontinued in December of 2021 as we couldn't find enough modelers to realistically finish the planned demo any time in the near future, much less the full game. We decided our time would be better invested in other projects, be they FusionFall related or not.</p><p>I worked on multiple things during my time on the team, such as movement, third person camera and generating animations for the main menu camera programmatically to save time whenever we'd make changes to it. But for the most part I was helping other volunteers by sharing informational resources from Unity's scripting documentation since I've learned my way around it relatively quickly.</p>"
$ws.Range("G4").Value = "https://fusionfall.fandom.com/wiki/FusionFall_RE:SPAWN"
$ws.Range("H4").Value = 0

# Register the hyperlink relationship for G4 (mirrors G2/G3). Adding a
# hyperlink via the collection re-applies Excel's built-in "Hyperlink" named
# style to the cell/font tables, so immediately restore the intended
# (already-matching) visual formatting and drop the now-unused named style.
$ws.Hyperlinks.Add($ws.Range("G4"), "https://fusionfall.fandom.com/wiki/FusionFall_RE:SPAWN") | Out-Null
Copy-FormatOnly "G3" "G4"
$wb.Styles("Hyperlink").Delete()
